$d = $word.ActiveDocument

# Find the paragraph with "Testing on 30th April 2023" (the last dated
# "Testing on ..." entry in the body) and its index within the
# Paragraphs collection.
$target = $null
$targetIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -match "Testing on 30th April 2023") {
        $target = $p
        $targetIndex = $i
    }
}

# Split the document right after that paragraph, creating a new
# (initially empty) paragraph that inherits the BodyText style.
$target.Range.InsertParagraphAfter()

# Fill in the freshly-created paragraph with the new sentence, using
# InsertXML so the run carries xml:space="preserve" like its siblings.
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$null = $newPara.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Testing on 19th May 2023</w:t></w:r></w:p>')
